$d = $word.ActiveDocument

$replacements = @(
    @("93×20=1860", "52×93=4836"),
    @("11×41=451", "46×11=506"),
    @("49×79=3871", "20×59=1180"),
    @("40×76=3040", "52×20=1040"),
    @("21×48=1008", "28×82=2296"),
    @("99×88=8712", "97×13=1261"),
    @("42×72=3024", "76×65=4940"),
    @("19×37=703", "24×68=1632"),
    @("32×16=512", "87×65=5655"),
    @("29×19=551", "33×27=891"),
    @("81×86=6966", "55×45=2475"),
    @("79×75=5925", "83×35=2905"),
    @("88×51=4488", "74×48=3552"),
    @("79×29=2291", "15×25=375"),
    @("58×58=3364", "78×20=1560"),
    @("47×28=1316", "64×43=2752"),
    @("79×76=6004", "48×86=4128"),
    @("81×20=1620", "18×87=1566"),
    @("38×57=2166", "11×31=341"),
    @("77×63=4851", "46×49=2254"),
    @("50×42=2100", "20×63=1260"),
    @("54×92=4968", "81×88=7128"),
    @("77×34=2618", "82×86=7052"),
    @("27×62=1674", "47×36=1692"),
    @("81×65=5265", "21×40=840")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
